$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: A10 no longer holds the stray 0 value ---
$ws.Range("A10").ClearContents()

# --- Column C: "mas_meh" -> "code_meh", values reformatted from "{x,y}" to "x, y" ---
$ws.Range("C1").Value2 = "code_meh"
$ws.Range("C2").Value2 = "6, 8"
$ws.Range("C3").ClearContents()
$ws.Range("C4").Value2 = "5, 6, 8"
$ws.Range("C5").Value2 = "5, 6, 8"
$ws.Range("C6").Value2 = "4, 5, 6, 8, 9"
$ws.Range("C7").Value2 = "3, 4, 5, 6, 8"
$ws.Range("C8").Value2 = "3, 5, 6, 9"
$ws.Range("C9").Value2 = "3, 4"

# --- Right-align the code_meh column (header + data cells) ---
foreach ($addr in @("C1","C2","C4","C5","C6","C7","C8","C9")) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Arial"
    $c.Font.Size = 10
    $c.HorizontalAlignment = -4152
}
$ws.Range("C3").HorizontalAlignment = -4152

# --- Widen name_poch column so the Russian soil names are fully visible ---
$ws.Columns("B").ColumnWidth = 24.6

# --- Leave the selection on C9, matching the saved view state ---
$ws.Range("C9").Select()
